$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert 3 new blank rows right before row 44. This pushes the
# trailing "Requerimiento" block (previously rows 44-52) down to rows
# 47-55, leaving rows 41-46 (3 previously-blank rows + 3 new ones) free
# for the new "userTypes" CRUD endpoints, matching the target layout
# (final used range A1:L55).
# ------------------------------------------------------------------
$ws.Rows("44:46").Insert()

# ------------------------------------------------------------------
# Fill in the new "userTypes" permit rows (41-46)
# ------------------------------------------------------------------
$ws.Range("B41").Value = "userTypes"
$ws.Range("C41").Value = "Admin"
$ws.Range("D41").Value = "GET"
$ws.Range("F41").Value = "Get all userTypes"

$ws.Range("B42").Value = "userTypes"
$ws.Range("C42").Value = "Admin"
$ws.Range("D42").Value = "GET"
$ws.Range("F42").Value = "Get Usertype by Id"

$ws.Range("B43").Value = "userTypes"
$ws.Range("C43").Value = "Admin"
$ws.Range("D43").Value = "POST"
$ws.Range("F43").Value = "Create new userType"

$ws.Range("B44").Value = "userTypes"
$ws.Range("C44").Value = "Admin"
$ws.Range("D44").Value = "PUT"
$ws.Range("F44").Value = "UPDATE userTypes"

$ws.Range("B45").Value = "userTypes"
$ws.Range("C45").Value = "Admin"
$ws.Range("D45").Value = "PATCH"
$ws.Range("F45").Value = "UPDATE userTypes"

$ws.Range("B46").Value = "userTypes"
$ws.Range("C46").Value = "Admin"
$ws.Range("D46").Value = "DELETE"
$ws.Range("F46").Value = "DELETE UserType"

# ------------------------------------------------------------------
# Cosmetic touch-ups matching the rest of the authored diff: a
# slightly narrower column B, a new (wider) column C, and the
# selection / scroll position left where the author ended up editing.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 22.27
$ws.Columns.Item(3).ColumnWidth = 19.6

$ws.Range("J44").Select()
